$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 12: average of J (the |S*|/n / k column) across the 10 result rows
$ws.Range("J12").Formula = "=AVERAGE(J2:J11)"
$ws.Range("J12").Font.Bold = $true

# Rows 14-17: summary labels + aggregate formulas
$ws.Range("A14").Value = "Average of SW(S*)/SW(OPT)"
$ws.Range("B14").Formula = "=AVERAGE(N2:N11)"

$ws.Range("A15").Value = "Average of SC(S*)/SC(OPT)"
$ws.Range("B15").Formula = "=AVERAGE(Z2:Z11)"

$ws.Range("A16").Value = "Worst of SW(S*)/SW(OPT)"
$ws.Range("B16").Formula = "=MIN(N2:N11)"

$ws.Range("A17").Value = "Worst of SC(S*)/SC(OPT)"
$ws.Range("B17").Formula = "=MAX(Z2:Z11)"

# Style the aggregate values in column B: bold, size 12, vertically centered
$rb = $ws.Range("B14:B17")
$rb.Font.Bold = $true
$rb.Font.Size = 12
$rb.VerticalAlignment = -4108

# Taller rows for the summary block
$ws.Range("A14:B17").RowHeight = 15.6

# Match the selection left behind in the saved file
$ws.Range("A14:B17").Select()

# Page setup tweaks present in the target file
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
